$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'25.912.67"
$ws.Cells.Item(2, 5).Value = "  +0.20%  "

$ws.Cells.Item(3, 4).Value = "'1.735.22"
$ws.Cells.Item(3, 5).Value = "  -0.25%  "

$ws.Cells.Item(4, 4).Value = "'0.9993"
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$ws.Cells.Item(5, 4).Value = "'246.01"
$ws.Cells.Item(5, 5).Value = "  +3.44%  "

$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 5).Value = "  +0.05%  "

$ws.Cells.Item(7, 4).Value = "'0.5018"
$ws.Cells.Item(7, 5).Value = "  -2.58%  "

$ws.Cells.Item(8, 4).Value = "'0.2719"
$ws.Cells.Item(8, 5).Value = "  -0.53%  "

$ws.Cells.Item(9, 4).Value = "'0.06168"
$ws.Cells.Item(9, 5).Value = "  +0.75%  "

$ws.Cells.Item(10, 2).Value = "TRON"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(10, 4).Value = "'0.07253"
$ws.Cells.Item(10, 5).Value = "  +1.10%  "

$ws.Cells.Item(11, 2).Value = "WrappedEther"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(11, 4).Value = "'1.736.32"
$ws.Cells.Item(11, 5).Value = "  -0.21%  "

$ws.Cells.Item(12, 4).Value = "'0.6532"
$ws.Cells.Item(12, 5).Value = "  +1.62%  "

$ws.Cells.Item(13, 4).Value = "'15.14"
$ws.Cells.Item(13, 5).Value = "  +1.46%  "

$ws.Cells.Item(14, 4).Value = "'4.784"
$ws.Cells.Item(14, 5).Value = "  +4.14%  "

$ws.Cells.Item(15, 4).Value = "'77.05"
$ws.Cells.Item(15, 5).Value = "  -0.33%  "

$ws.Cells.Item(16, 5).Value = "  +0.14%  "

$ws.Cells.Item(17, 4).Value = "'0.9990"
$ws.Cells.Item(17, 5).Value = "  -0.06%  "

$ws.Cells.Item(18, 4).Value = "'25.917.76"
$ws.Cells.Item(18, 5).Value = "  +0.17%  "

$ws.Cells.Item(19, 4).Value = "'11.83"
$ws.Cells.Item(19, 5).Value = "  +0.85%  "

$ws.Cells.Item(20, 4).Value = "'0.000006806"
$ws.Cells.Item(20, 5).Value = "  +0.63%  "

$ws.Cells.Item(21, 4).Value = "'4.610"
$ws.Cells.Item(21, 5).Value = "  +8.26%  "

$ws.Cells.Item(22, 4).Value = "'1.958.82"
$ws.Cells.Item(22, 5).Value = "  -0.21%  "

$ws.Cells.Item(23, 4).Value = "'8.768"
$ws.Cells.Item(23, 5).Value = "  +1.11%  "

$ws.Cells.Item(24, 4).Value = "'5.477"
$ws.Cells.Item(24, 5).Value = "  +4.57%  "

$ws.Cells.Item(25, 4).Value = "'133.79"
$ws.Cells.Item(25, 5).Value = "  -3.62%  "

$ws.Cells.Item(26, 4).Value = "'15.26"
$ws.Cells.Item(26, 5).Value = "  +0.63%  "

$ws.Cells.Item(27, 4).Value = "'1.785"
$ws.Cells.Item(27, 5).Value = "  +1.40%  "

$ws.Cells.Item(28, 4).Value = "'1.412"
$ws.Cells.Item(28, 5).Value = "  -6.74%  "

$ws.Cells.Item(29, 4).Value = "'105.63"
$ws.Cells.Item(29, 5).Value = "  -0.17%  "

$ws.Cells.Item(30, 4).Value = "'3.996"
$ws.Cells.Item(30, 5).Value = "  -0.63%  "

$ws.Cells.Item(31, 4).Value = "'0.08108"
$ws.Cells.Item(31, 5).Value = "  -2.45%  "

$ws.Cells.Item(32, 4).Value = "'3.709"
$ws.Cells.Item(32, 5).Value = "  +1.91%  "

$ws.Cells.Item(33, 4).Value = "'0.04731"
$ws.Cells.Item(33, 5).Value = "  +3.03%  "

$ws.Cells.Item(34, 4).Value = "'2.657"
$ws.Cells.Item(34, 5).Value = "  -0.21%  "

$ws.Cells.Item(35, 4).Value = "'0.9971"
$ws.Cells.Item(35, 5).Value = "  +0.89%  "

$ws.Cells.Item(36, 4).Value = "'0.6110"
$ws.Cells.Item(36, 5).Value = "  -1.26%  "

$ws.Cells.Item(37, 4).Value = "'2.736"
$ws.Cells.Item(37, 5).Value = "  +1.90%  "

$ws.Cells.Item(38, 4).Value = "'0.01602"
$ws.Cells.Item(38, 5).Value = "  -0.77%  "

$ws.Cells.Item(39, 4).Value = "'0.8549"
$ws.Cells.Item(39, 5).Value = "  +16.22%  "

$ws.Cells.Item(40, 4).Value = "'1.939"
$ws.Cells.Item(40, 5).Value = "  +0.44%  "

$ws.Cells.Item(41, 4).Value = "'1.0000"
$ws.Cells.Item(41, 5).Value = "  +0.05%  "

$ws.Cells.Item(42, 4).Value = "'100.51"
$ws.Cells.Item(42, 5).Value = "  +2.89%  "

$ws.Cells.Item(43, 4).Value = "'0.3911"
$ws.Cells.Item(43, 5).Value = "  +1.88%  "

$ws.Cells.Item(44, 4).Value = "'5.016"
$ws.Cells.Item(44, 5).Value = "  +1.25%  "

$ws.Cells.Item(45, 4).Value = "'0.1176"
$ws.Cells.Item(45, 5).Value = "  +4.59%  "

$ws.Cells.Item(46, 4).Value = "'6.327"
$ws.Cells.Item(46, 5).Value = "  +2.52%  "

$ws.Cells.Item(47, 4).Value = "'55.76"
$ws.Cells.Item(47, 5).Value = "  +1.64%  "

$ws.Cells.Item(48, 4).Value = "'0.05273"

$ws.Cells.Item(49, 4).Value = "'30.71"
$ws.Cells.Item(49, 5).Value = "  +0.76%  "

$ws.Cells.Item(50, 2).Value = "Decentraland"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(50, 4).Value = "'0.3475"
$ws.Cells.Item(50, 5).Value = "  +1.87%  "

$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "'7.598"
$ws.Cells.Item(51, 5).Value = "  -0.05%  "
